$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: remove K1 entirely (content+format) -> header no longer has a K column
$ws.Range("K1").Clear()

# Rows 2-26: the J cell was an empty placeholder and K held the real
# "didnotparticipate" value. Move that value into J, then clear K entirely.
for ($r = 2; $r -le 26; $r++) {
    $kVal = $ws.Cells.Item($r, 11).Value2
    $ws.Cells.Item($r, 10).Value = $kVal
    $ws.Cells.Item($r, 11).Clear()
}

# Rows 27-70: J already holds the real value; K is just an empty placeholder -> clear it.
for ($r = 27; $r -le 70; $r++) {
    $ws.Cells.Item($r, 11).Clear()
}

# Fix E53 trailing-space typo: "concurrence " -> "concurrence"
$ws.Range("E53").Value = "concurrence"
